$wb = $excel.ActiveWorkbook

# --- Rename sheets (task-order run ids refreshed) ---
$wb.Worksheets.Item(1).Name = "GNG_TO-1650996072255417"
$wb.Worksheets.Item(2).Name = "NB_TO-16509960739195604"
$wb.Worksheets.Item(3).Name = "RS_TO-16509960739195604"
$wb.Worksheets.Item(4).Name = "TOL_TO-16509960739755964"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16509960740635993"

# --- Sheet 1 (GNG) ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("B2").Value = "go_stims-16509960722233825.csv"
$ws.Range("B3").Value = "GNG_stims-16509960722393818.csv"
$ws.Range("B4").Value = "go_stims-16509960722393818.csv"
$ws.Range("B5").Value = "GNG_stims-1650996072255417.csv"

# --- Sheet 2 (NB) ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("B2").Value = "ZB-match_5-1650996072431426.csv"
$ws.Range("B3").Value = "TB-16509960734955962.csv"
$ws.Range("B4").Value = "TB-16509960731035974.csv"
$ws.Range("B5").Value = "OB-16509960729435945.csv"
$ws.Range("B6").Value = "TB-16509960738955686.csv"
$ws.Range("B7").Value = "ZB-match_7-1650996072271424.csv"
$ws.Range("B8").Value = "OB-16509960726715717.csv"
$ws.Range("B9").Value = "OB-16509960727195618.csv"
$ws.Range("B10").Value = "ZB-match_8-1650996072487383.csv"

# --- Sheet 4 (TOL) ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("B2").Value = "MM_stims-16509960739435663.csv"
$ws.Range("B3").Value = "ZM_stims-16509960739195604.csv"
$ws.Range("B4").Value = "MM_stims-16509960739595957.csv"
$ws.Range("B5").Value = "ZM_stims-16509960739435663.csv"
$ws.Range("B6").Value = "MM_stims-16509960739755964.csv"
$ws.Range("B7").Value = "ZM_stims-16509960739595957.csv"

# --- Sheet 5 (vSAT) ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("B2").Value = "SAT_stims-16509960739755964.csv"
$ws.Range("B3").Value = "vSAT_stims-16509960740235965.csv"
$ws.Range("B4").Value = "vSAT_stims-16509960740395606.csv"
$ws.Range("B5").Value = "SAT_stims-16509960740075984.csv"
